$wb = $excel.ActiveWorkbook

# --- Rename existing sheets ---
$wb.Worksheets.Item("mob").Name = "mobility"
$wb.Worksheets.Item("inf").Name = "infrastructure"
$wb.Worksheets.Item("deriv").Name = "derivation"

# --- Update the Print_Area defined name to follow the renamed sheet ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "mobility!Print_Area") {
        $n.RefersTo = "=mobility!`$A`$1:`$C`$20"
    }
}

# --- Add the new "categories" sheet after "derivation" ---
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item("derivation"))
$ws4.Name = "categories"

# Tab color (theme 7, tint -0.249977111117893 -> resolved RGB 60497A) to
# match the other sheets (the other tabs already resolve to this same RGB)
$ws4.Tab.Color = 8014176

# Column widths (approximate values closest to 21.140625 / 10.5703125 /
# 107.7109375 character-units that this engine's pixel-quantized
# ColumnWidth setter can represent)
$ws4.Columns.Item(1).ColumnWidth = 20.333333333333332
$ws4.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws4.Columns.Item(3).ColumnWidth = 106.83333333333333

# Header row
$ws4.Range("A1").Value = "Variable name"
$ws4.Range("B1").Value = "Value"
$ws4.Range("C1").Value = "Description"
$ws4.Range("A1:C1").Font.Bold = $true
$ws4.Range("B1").HorizontalAlignment = -4108

# Data rows
$names = @("max_derivation_1", "max_derivation_2", "max_derivation_3", "max_derivation_4", "max_derivation_5", "regroup_1", "regroup_2", "regroup_3", "regroup_4", "regroup_5")
$values = @(0.8, 0.7, 0.7, 0.6, 0.5, 1, 1, 0, 0, 0)
$descriptions = @(
    "Maximum % of derivation for product category 1 - grains (coeff).",
    "Maximum % of derivation for product category 2 - primary products no grains (coeff).",
    "Maximum % of derivation for product category 3 - semi manufactured (coeff).",
    "Maximum % of derivation for product category 4 - manufactured (coeff).",
    "Maximum % of derivation for product category 5 - unknown (coeff).",
    "Regroup trains to reduce idle capacity (1=yes, 0=no) for product category 1 - grains (coeff).",
    "Regroup trains to reduce idle capacity (1=yes, 0=no) for product category 2 - primary products no grains (coeff).",
    "Regroup trains to reduce idle capacity (1=yes, 0=no) for product category 3 - semi manufactured (coeff).",
    "Regroup trains to reduce idle capacity (1=yes, 0=no) for product category 4 - manufactured (coeff).",
    "Regroup trains to reduce idle capacity (1=yes, 0=no) for product category 5 - unknown (coeff)."
)

# Populate column A (variable names) for all rows first, then column C
# (descriptions), so new shared strings are appended in the same order
# as the target workbook (all names, then all descriptions).
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws4.Cells.Item($i + 2, 1).Value = $names[$i]
}
for ($i = 0; $i -lt $descriptions.Length; $i++) {
    $ws4.Cells.Item($i + 2, 3).Value = $descriptions[$i]
}
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws4.Cells.Item($i + 2, 2).Value = $values[$i]
}

$ws4.Range("C3").Select()

# --- Adjust selections on existing sheets to match final state ---
$ws1 = $wb.Worksheets.Item("mobility")
$ws1.Range("B20").Select()

$ws3 = $wb.Worksheets.Item("derivation")
$ws3.Range("B4").Select()

$ws4.Activate()
